$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 469. This shifts all rows 469..520
# down to 470..521 automatically (Excel row-insert semantics), leaving the
# new row 469 empty and ready to be populated.
$ws.Rows.Item(469).Insert()

# Populate the newly inserted row 469 with the record that used to occupy
# row 468 (date 44468 / "Sin especificar" / "2a amarillo" / ...).
$ws.Cells.Item(469, 1).Value = 1
$ws.Cells.Item(469, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(469, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(469, 4).Value = 44468
$ws.Cells.Item(469, 5).Value = 15
$ws.Cells.Item(469, 6).Value = "Fruta"
$ws.Cells.Item(469, 7).Value = 100102
$ws.Cells.Item(469, 8).Value = "Cítricos"
$ws.Cells.Item(469, 9).Value = 100102003
$ws.Cells.Item(469, 10).Value = "Limón"
$ws.Cells.Item(469, 11).Value = "Sin especificar"
$ws.Cells.Item(469, 12).Value = "2a amarillo"
$ws.Cells.Item(469, 13).Value = 200
$ws.Cells.Item(469, 14).Value = 11000
$ws.Cells.Item(469, 15).Value = 12000
$ws.Cells.Item(469, 16).Value = 11500
$ws.Cells.Item(469, 17).Value = '$/caja 20 kilos'
$ws.Cells.Item(469, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(469, 19).Value = 575
$ws.Cells.Item(469, 20).Value = 20

# Row 468 itself gets a new date and new price-grade figures (a fresh
# "1a amarillo" quotation replacing what had been entered there before).
$ws.Cells.Item(468, 4).Value = 45223
$ws.Cells.Item(468, 12).Value = "1a amarillo"
$ws.Cells.Item(468, 13).Value = 300
$ws.Cells.Item(468, 14).Value = 17000
$ws.Cells.Item(468, 15).Value = 18000
$ws.Cells.Item(468, 16).Value = 17500
$ws.Cells.Item(468, 19).Value = 875
